$d = $word.ActiveDocument

# Locate the paragraph that contains the astromap credit/link line.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Jan Hollan, CzechGlobe*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $full = $target.Range
    # Exclude the trailing paragraph mark from the range.
    $full.End = $full.End - 1
    $start = $full.Start
    $oldLen = $full.End - $full.Start

    $newText = "de Jan Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

    # Build a minimal OOXML fragment: a leading empty run followed by a
    # single unformatted run carrying the whole merged text, matching the
    # target structure exactly (bypasses run auto-merge/format inheritance
    # that plain Range.Text assignment would otherwise apply).
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    # Insert the replacement content immediately before the old text.
    $ins = $d.Range($start, $start)
    $ins.InsertXML($xmlFrag)

    # Remove the old run text, now shifted past the freshly inserted text.
    $newLen = $newText.Length
    $oldRange = $d.Range($start + $newLen, $start + $newLen + $oldLen)
    $oldRange.Text = ""
}
